$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 4 new event rows, working from the top of the sheet downward.
# Each Insert() happens at the row's final target position; since all
# earlier (lower-numbered) insertions have already been applied, the rows
# above are already settled in their final place.

# New row -> final position 6 ("Illusionary Line Art Workshop", before "Jhumritalaiya")
$ws.Rows("6").Insert()
$ws.Range("A6").Value = "Illusionary Line Art Workshop"
$ws.Range("B6").Value = "December 19 | 3PM"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "150"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "Online"
$ws.Range("E6").Value = "online-hobby"

# New row -> final position 8 ("Learn Mandala Art ...", before "Lose 3 Kgs")
$ws.Rows("8").Insert()
$ws.Range("A8").Value = "Learn Mandala Art - A One Page Spotlight Workshop"
$ws.Range("B8").Value = "September 27 | 12PM - December 31 | 6PM"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "699"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "Online"
$ws.Range("E8").Value = "online-hobby"

# New row -> final position 12 ("Online Internship Program ...", before "Online Program in Direction")
$ws.Rows("12").Insert()
$ws.Range("A12").Value = "Online Internship Program for Drone/Robotics/IoT/Computer Vision"
$ws.Range("B12").Value = "August 21 | 10AM - December 31 | 4PM"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "1000"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "To Be Announced"
$ws.Range("E12").Value = "online-hobby"

# New row -> final position 16 ("Paper Mache ...", before SHARAN)
$ws.Rows("16").Insert()
$ws.Range("A16").Value = "Paper Mache [Online Live Workshop - Inclusive of Materials]"
$ws.Range("B16").Value = "December 19 | 11AM"
$ws.Range("C16").Value = "1500 onwards"
$ws.Range("D16").Value = "Zoom Call"
$ws.Range("E16").Value = "online-hobby"
